$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from SCD0174 to SCD0011
$ws.Name = "SCD0011"

# Update TC_ID values in B2/B3 from "DGS-189" to "SCD0011-005"
$ws.Range("B2").Value = "SCD0011-005"
$ws.Range("B3").Value = "SCD0011-005"

# Widen column B to fit the new, longer TC_ID text
$ws.Columns.Item(2).ColumnWidth = 11.67

# Move the active selection to B4
$ws.Range("B4").Select()
